# chore: update Sheets via scheduled runner
# Refresh cached Market Board price columns (H/I/J/K/L) and the
# derived profit columns (M/N) for the affected Leve rows on each
# job sheet, reflecting the latest price pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 410.21054
$ws.Range("I15").Value = 410.21054
$ws.Range("K15").Value = 1230.63162
$ws.Range("M15").Value = -1061.63162
$ws.Range("H80").Value = 1709.3334
$ws.Range("J80").Value = 2081.2
$ws.Range("L80").Value = 6243.599999999999
$ws.Range("N80").Value = -8239.599999999999
$ws.Range("H83").Value = 1709.3334
$ws.Range("J83").Value = 2081.2
$ws.Range("L83").Value = 18730.8
$ws.Range("N83").Value = -28714.8
$ws.Range("H92").Value = 313.1111
$ws.Range("I92").Value = 227.375
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 227.375
$ws.Range("L92").Value = 999
$ws.Range("M92").Value = 1020.625
$ws.Range("N92").Value = -3495
$ws.Range("H96").Value = 500
$ws.Range("J96").Value = 500
$ws.Range("L96").Value = 1500
$ws.Range("N96").Value = -4246
$ws.Range("H97").Value = 2000
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H104").Value = 2866.3333
$ws.Range("I104").Value = 2866.3333
$ws.Range("K104").Value = 8598.999899999999
$ws.Range("M104").Value = -6851.999899999999
$ws.Range("H112").Value = 1200
$ws.Range("J112").Value = 1200
$ws.Range("L112").Value = 3600
$ws.Range("N112").Value = -5816
$ws.Range("H118").Value = 832.25
$ws.Range("I118").Value = 776.3333
$ws.Range("K118").Value = 2328.9999
$ws.Range("M118").Value = -671.9998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 363.22223
$ws.Range("I110").Value = 291.25
$ws.Range("J110").Value = 420.8
$ws.Range("K110").Value = 291.25
$ws.Range("L110").Value = 420.8
$ws.Range("M110").Value = 1753.75
$ws.Range("N110").Value = -4510.8
$ws.Range("H122").Value = 7615.3335
$ws.Range("I122").Value = 8896.4
$ws.Range("K122").Value = 26689.2
$ws.Range("M122").Value = -24239.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 612.5
$ws.Range("I94").Value = 612.5
$ws.Range("K94").Value = 612.5
$ws.Range("M94").Value = -161.5
$ws.Range("H107").Value = 1798
$ws.Range("I107").Value = 899.5
$ws.Range("K107").Value = 899.5
$ws.Range("M107").Value = 1020.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2131.5
$ws.Range("I35").Value = 1508.6666
$ws.Range("J35").Value = 4000
$ws.Range("K35").Value = 1508.6666
$ws.Range("L35").Value = 4000
$ws.Range("M35").Value = -1214.6666
$ws.Range("N35").Value = -4588
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 39983.332
$ws.Range("J74").Value = 39983.332
$ws.Range("L74").Value = 39983.332
$ws.Range("N74").Value = -41731.332
$ws.Range("H77").Value = 39983.332
$ws.Range("J77").Value = 39983.332
$ws.Range("L77").Value = 119949.996
$ws.Range("N77").Value = -128685.996
$ws.Range("H107").Value = 392.83334
$ws.Range("I107").Value = 392.83334
$ws.Range("K107").Value = 392.83334
$ws.Range("M107").Value = 1527.16666
$ws.Range("H134").Value = 5444.909
$ws.Range("I134").Value = 2842.2856
$ws.Range("K134").Value = 8526.856800000001
$ws.Range("M134").Value = -5991.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 462.5
$ws.Range("J4").Value = 50
$ws.Range("L4").Value = 150
$ws.Range("N4").Value = -374
$ws.Range("H13").Value = 273
$ws.Range("I13").Value = 273
$ws.Range("K13").Value = 819
$ws.Range("M13").Value = -651
$ws.Range("H61").Value = 59
$ws.Range("I61").Value = 31.666666
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 94.99999800000001
$ws.Range("L61").Value = 300
$ws.Range("M61").Value = 120.000002
$ws.Range("N61").Value = -730
$ws.Range("H94").Value = 2535.625
$ws.Range("I94").Value = 3962
$ws.Range("J94").Value = 2060.1667
$ws.Range("K94").Value = 11886
$ws.Range("L94").Value = 6180.500100000001
$ws.Range("M94").Value = -11210
$ws.Range("N94").Value = -7532.500100000001
$ws.Range("H113").Value = 378.66666
$ws.Range("I113").Value = 384
$ws.Range("J113").Value = 376
$ws.Range("K113").Value = 1152
$ws.Range("L113").Value = 1128
$ws.Range("M113").Value = 1018
$ws.Range("N113").Value = -5468
$ws.Range("H115").Value = 1411.25
$ws.Range("I115").Value = 1659.6666
$ws.Range("J115").Value = 666
$ws.Range("K115").Value = 4978.9998
$ws.Range("L115").Value = 1998
$ws.Range("M115").Value = -3803.9998
$ws.Range("N115").Value = -4348

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9666
$ws.Range("I7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("M7").Value = -1887
$ws.Range("H22").Value = 2339.4443
$ws.Range("I22").Value = 938.75
$ws.Range("J22").Value = 3460
$ws.Range("K22").Value = 938.75
$ws.Range("L22").Value = 3460
$ws.Range("M22").Value = -643.75
$ws.Range("N22").Value = -4050
$ws.Range("H27").Value = 2339.4443
$ws.Range("I27").Value = 938.75
$ws.Range("J27").Value = 3460
$ws.Range("K27").Value = 938.75
$ws.Range("L27").Value = 3460
$ws.Range("M27").Value = -831.75
$ws.Range("N27").Value = -3674
$ws.Range("H122").Value = 4008.6
$ws.Range("I122").Value = 3810.875
$ws.Range("K122").Value = 11432.625
$ws.Range("M122").Value = -8982.625
$ws.Range("H126").Value = 9666
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527
$ws.Range("H137").Value = 47100
$ws.Range("I137").Value = 47100
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 47100
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -42000
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3606.2856
$ws.Range("I107").Value = 1411.6666
$ws.Range("J107").Value = 4234.9998
$ws.Range("K107").Value = 4234.9998
$ws.Range("M107").Value = -2314.9998
